$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 11 (Leve Item ID 5533)
$ws.Range("H11").Value = 678.35297
$ws.Range("I11").Value = 678.35297
$ws.Range("K11").Value = 678.35297
$ws.Range("M11").Value = -538.35297

# Row 17 (Leve Item ID 38956)
$ws.Range("H17").Value = 1186.3334
$ws.Range("J17").Value = 1186.3334
$ws.Range("L17").Value = 3559.0002
$ws.Range("N17").Value = -3895.0002

# Row 18 (Leve Item ID 5471)
$ws.Range("H18").Value = 8460.200000000001
$ws.Range("I18").Value = 8460.200000000001
$ws.Range("K18").Value = 8460.200000000001
$ws.Range("M18").Value = -8176.200000000001

# Row 19 (Leve Item ID 7015)
$ws.Range("H19").Value = 245
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 245
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 245
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -595

# Row 62 (Leve Item ID 27781)
$ws.Range("H62").Value = 9938.615
$ws.Range("I62").Value = 5424.5
$ws.Range("J62").Value = 11944.889
$ws.Range("K62").Value = 5424.5
$ws.Range("L62").Value = 11944.889
$ws.Range("M62").Value = -4800.5
$ws.Range("N62").Value = -13192.889

# Row 65 (Leve Item ID 27781)
$ws.Range("H65").Value = 9938.615
$ws.Range("I65").Value = 5424.5
$ws.Range("J65").Value = 11944.889
$ws.Range("K65").Value = 27122.5
$ws.Range("L65").Value = 59724.44499999999
$ws.Range("M65").Value = -24002.5
$ws.Range("N65").Value = -65964.44499999999

# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 5799.6
$ws.Range("I116").Value = 8233
$ws.Range("J116").Value = 2149.5
$ws.Range("K116").Value = 8233
$ws.Range("L116").Value = 2149.5
$ws.Range("M116").Value = -4791
$ws.Range("N116").Value = -9033.5

# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 2814
$ws.Range("I132").Value = 2014.9445
$ws.Range("K132").Value = 6044.833500000001
$ws.Range("M132").Value = -3514.833500000001

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 9369.117
$ws.Range("I137").Value = 12589
$ws.Range("K137").Value = 37767
$ws.Range("M137").Value = -35217

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 3462.2856
$ws.Range("I138").Value = 2732.3845
$ws.Range("J138").Value = 3682.9534
$ws.Range("K138").Value = 8197.1535
$ws.Range("L138").Value = 11048.8602
$ws.Range("M138").Value = -3057.1535
$ws.Range("N138").Value = -21328.8602


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 3335.907
$ws.Range("I32").Value = 3492.795
$ws.Range("K32").Value = 3492.795
$ws.Range("M32").Value = -3205.795

# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 3644.7
$ws.Range("J45").Value = 4378.143
$ws.Range("L45").Value = 4378.143
$ws.Range("N45").Value = -5132.143

# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 6308.643
$ws.Range("I61").Value = 6555.4614
$ws.Range("K61").Value = 6555.4614
$ws.Range("M61").Value = -6343.4614

# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 3579.15
$ws.Range("I74").Value = 3613.0715
$ws.Range("K74").Value = 3613.0715
$ws.Range("M74").Value = -2739.0715

# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 3579.15
$ws.Range("I77").Value = 3613.0715
$ws.Range("K77").Value = 18065.3575
$ws.Range("M77").Value = -13697.3575

# Row 97 (Leve Item ID 19941)
$ws.Range("H97").Value = 620.2143
$ws.Range("I97").Value = 144.8
$ws.Range("K97").Value = 144.8
$ws.Range("M97").Value = 351.2

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 3127.818
$ws.Range("I132").Value = 3312.7273
$ws.Range("J132").Value = 2942.9092
$ws.Range("K132").Value = 9938.1819
$ws.Range("L132").Value = 8828.7276
$ws.Range("M132").Value = -7408.1819
$ws.Range("N132").Value = -13888.7276

# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 6308.643
$ws.Range("I136").Value = 6555.4614
$ws.Range("K136").Value = 19666.3842
$ws.Range("M136").Value = -17116.3842


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 145447.05
$ws.Range("I94").Value = 166096.72
$ws.Range("J94").Value = 899.3333
$ws.Range("K94").Value = 166096.72
$ws.Range("L94").Value = 899.3333
$ws.Range("M94").Value = -165645.72
$ws.Range("N94").Value = -1801.3333

# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 6666.6665
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 6666.6665
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 6666.6665
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -9662.666499999999

# Row 109 (Leve Item ID 27096)
$ws.Range("H109").Value = 67981
$ws.Range("J109").Value = 67981
$ws.Range("L109").Value = 67981
$ws.Range("N109").Value = -70755

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 3838.875
$ws.Range("I134").Value = 4528
$ws.Range("J134").Value = 3149.75
$ws.Range("K134").Value = 13584
$ws.Range("L134").Value = 9449.25
$ws.Range("M134").Value = -11049
$ws.Range("N134").Value = -14519.25


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 87 (Leve Item ID 11929)
$ws.Range("H87").Value = 59999.5
$ws.Range("J87").Value = 59999.5
$ws.Range("L87").Value = 59999.5
$ws.Range("N87").Value = -62371.5

# Row 90 (Leve Item ID 11929)
$ws.Range("H90").Value = 59999.5
$ws.Range("J90").Value = 59999.5
$ws.Range("L90").Value = 179998.5
$ws.Range("N90").Value = -191854.5

# Row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 1998.0625
$ws.Range("I99").Value = 1510.375
$ws.Range("J99").Value = 2485.75
$ws.Range("K99").Value = 1510.375
$ws.Range("L99").Value = 2485.75
$ws.Range("M99").Value = -12.375
$ws.Range("N99").Value = -5481.75

# Row 105 (Leve Item ID 19928)
$ws.Range("H105").Value = 3046.0908
$ws.Range("I105").Value = 2251.3333
$ws.Range("J105").Value = 3999.8
$ws.Range("K105").Value = 2251.3333
$ws.Range("L105").Value = 3999.8
$ws.Range("M105").Value = -504.3332999999998
$ws.Range("N105").Value = -7493.8

# Row 122 (Leve Item ID 36196)
$ws.Range("H122").Value = 1961.7931
$ws.Range("I122").Value = 1617.7142
$ws.Range("J122").Value = 2865
$ws.Range("K122").Value = 4853.142599999999
$ws.Range("L122").Value = 8595
$ws.Range("M122").Value = -2403.142599999999
$ws.Range("N122").Value = -13495

# Row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 1998.0625
$ws.Range("I126").Value = 1510.375
$ws.Range("J126").Value = 2485.75
$ws.Range("K126").Value = 4531.125
$ws.Range("L126").Value = 7457.25
$ws.Range("M126").Value = -2061.125
$ws.Range("N126").Value = -12397.25

# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 1764.3334
$ws.Range("I132").Value = 1764.3334
$ws.Range("K132").Value = 5293.0002
$ws.Range("M132").Value = -2763.0002


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 46 (Leve Item ID 4701)
$ws.Range("H46").Value = 14949986
$ws.Range("I46").Value = 13939983
$ws.Range("J46").Value = 20000000
$ws.Range("K46").Value = 41819949
$ws.Range("L46").Value = 60000000
$ws.Range("M46").Value = -41819858
$ws.Range("N46").Value = -60000182

# Row 104 (Leve Item ID 19807)
$ws.Range("H104").Value = 14536.223
$ws.Range("I104").Value = 8332.333000000001
$ws.Range("J104").Value = 17638.166
$ws.Range("K104").Value = 24996.999
$ws.Range("L104").Value = 52914.49800000001
$ws.Range("M104").Value = -22375.999
$ws.Range("N104").Value = -58156.49800000001


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 3319.9443
$ws.Range("I122").Value = 3221.1177
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 9663.3531
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -7213.3531
$ws.Range("N122").Value = -19900


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 3193.2632
$ws.Range("I22").Value = 2627.7646
$ws.Range("J22").Value = 8000
$ws.Range("K22").Value = 2627.7646
$ws.Range("L22").Value = 8000
$ws.Range("M22").Value = -2332.7646
$ws.Range("N22").Value = -8590

# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 3193.2632
$ws.Range("I27").Value = 2627.7646
$ws.Range("J27").Value = 8000
$ws.Range("K27").Value = 2627.7646
$ws.Range("L27").Value = 8000
$ws.Range("M27").Value = -2520.7646
$ws.Range("N27").Value = -8214

# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 2871
$ws.Range("I46").Value = 919.1429000000001
$ws.Range("J46").Value = 3922
$ws.Range("K46").Value = 919.1429000000001
$ws.Range("L46").Value = 3922
$ws.Range("M46").Value = -731.1429000000001
$ws.Range("N46").Value = -4298

# Row 61 (Leve Item ID 27740)
$ws.Range("H61").Value = 5113.65
$ws.Range("I61").Value = 2186.111
$ws.Range("K61").Value = 2186.111
$ws.Range("M61").Value = -1984.111

# Row 82 (Leve Item ID 12565)
$ws.Range("H82").Value = 428194.38
$ws.Range("I82").Value = 757145.75
$ws.Range("K82").Value = 757145.75
$ws.Range("M82").Value = -756784.75

# Row 85 (Leve Item ID 12565)
$ws.Range("H85").Value = 428194.38
$ws.Range("I85").Value = 757145.75
$ws.Range("K85").Value = 757145.75
$ws.Range("M85").Value = -755897.75

# Row 113 (Leve Item ID 27740)
$ws.Range("H113").Value = 5113.65
$ws.Range("I113").Value = 2186.111
$ws.Range("K113").Value = 2186.111
$ws.Range("M113").Value = -16.11099999999988

# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 3138.4348
$ws.Range("I136").Value = 3191
$ws.Range("K136").Value = 9573
$ws.Range("M136").Value = -7023


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 2 (Leve Item ID 3307)
$ws.Range("H2").Value = 371809.7
$ws.Range("I2").Value = 412010.88
$ws.Range("J2").Value = 9999
$ws.Range("K2").Value = 412010.88
$ws.Range("L2").Value = 9999
$ws.Range("M2").Value = -411898.88
$ws.Range("N2").Value = -10223

# Row 45 (Leve Item ID 21726)
$ws.Range("H45").Value = 9627.333000000001
$ws.Range("I45").Value = 9500
$ws.Range("K45").Value = 9500
$ws.Range("M45").Value = -9009

# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 2409.353
$ws.Range("I132").Value = 2320.258
$ws.Range("J132").Value = 3330
$ws.Range("K132").Value = 6960.773999999999
$ws.Range("L132").Value = 9990
$ws.Range("M132").Value = -4430.773999999999
$ws.Range("N132").Value = -15050

# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 2518.8572
$ws.Range("I136").Value = 2284
$ws.Range("K136").Value = 6852
$ws.Range("M136").Value = -4302

